$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 2D training schedule values (rows 2-6, columns B:H)
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = -1
$ws.Range("H2").Value = 65

$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 4

$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 7
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = -3
$ws.Range("H4").Value = 43

$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = -5
$ws.Range("H5").Value = 21

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 5

# Update selection so the sheet view records the active cell as I1
$ws.Range("I1").Select()
